# Apply cryptos list update (price & volume changes) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.123.02"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "2.018.01"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'226.64"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").Value = "'0.607"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'54.94"
$ws.Range("E8").Value = "  -3.88%  "
$ws.Range("D9").Value = "'0.380"
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "'0.0788"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  -3.18%  "
$ws.Range("D12").Value = "2.318.16"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("D13").Value = "'14.26"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("D14").Value = "'20.45"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "'0.742"
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").Value = "'5.14"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").Value = "1.995.64"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("D18").Value = "37.027.44"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").Value = "'6.16"
$ws.Range("E19").Value = "  +2.78%  "
$ws.Range("D20").Value = "'68.90"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("D21").Value = "0.0₃0822"
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "'224.47"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("E25").Value = "  -5.02%  "
$ws.Range("D26").Value = "'165.51"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("D27").Value = "'9.19"
$ws.Range("E27").Value = "  -3.20%  "
$ws.Range("D28").Value = "'0.126"
$ws.Range("E28").Value = "  -3.92%  "
$ws.Range("D29").Value = "'1.37"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").Value = "'18.74"
$ws.Range("E30").Value = "  -2.23%  "
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").Value = "'4.56"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'0.0618"
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("D34").Value = "'4.42"
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("D35").Value = "'2.35"
$ws.Range("E35").Value = "  -5.31%  "
$ws.Range("E36").Value = "  +1.68%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "'3.15"
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("D39").Value = "'5.41"
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.485.34"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0217"
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").Value = "'95.33"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("D43").Value = "'16.62"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").Value = "'0.0925"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("E45").Value = "  -4.73%  "
$ws.Range("E46").Value = "  -4.49%  "
$ws.Range("D47").Value = "'7.32"
$ws.Range("E48").Value = "  -2.17%  "
$ws.Range("D49").Value = "'2.91"
$ws.Range("E49").Value = "  -0.68%  "
$ws.Range("D50").Value = "2.208.05"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").Value = "'44.37"
$ws.Range("E51").Value = "  -1.90%  "
